# Append the new daily profit row (row 57, dated 01/20/2026) to Sheet1,
# matching the existing table's column layout:
#   A: Date (text)            B: Portfolio Value(USD)
#   C: BTC %                  D: KAS %
#   E: KAS Profit(USD)        F: KAS Profit(%)
#   G: KAS Total Profit(USD)  H: KAS Total Profit(%)
#   I: BTC Profit(USD)        J: BTC Profit(%)
#   K: Combined Total Profit(USD)  L: Combined Total Profit(%)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 57

# Column A holds the date as plain text in every existing row (e.g.
# "01/19/2026" in A56), not a real Excel date serial. A leading apostrophe
# forces Excel to store the literal text instead of autoconverting the
# "MM/DD/YYYY" string into a date value; resetting the style back to
# "Normal" afterwards drops the resulting quote-prefix formatting so the
# new cell ends up unstyled, just like the other data rows.
$ws.Cells.Item($row, 1).Value = "'01/20/2026"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 11973.81
$ws.Cells.Item($row, 3).Value = 0.2355691143938523
$ws.Cells.Item($row, 4).Value = 0.7644308856061477
$ws.Cells.Item($row, 5).Value = -178.59
$ws.Cells.Item($row, 6).Value = -26.07
$ws.Cells.Item($row, 7).Value = -21588.6
$ws.Cells.Item($row, 8).Value = -70.23
$ws.Cells.Item($row, 9).Value = -382.19
$ws.Cells.Item($row, 10).Value = -11.93
$ws.Cells.Item($row, 11).Value = -21970.79
$ws.Cells.Item($row, 12).Value = -64.73
